$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Merge the "Versi" + "on" runs back into a single "Version" run (no text change,
# but forces Word to rebuild the run across the old run boundary).
$find.Execute("Version", $true, $false, $false, $false, $false, $true, 1, $false, "Version", 2)

# Update the version number run " 2" -> " 1." (adds the trailing period here).
$find.Execute(" 2", $true, $false, $false, $false, $false, $true, 1, $false, " 1.", 2)

# Remove the now-redundant trailing "." run that followed the bookmark.
$find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
